$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy formatting from the
# existing header cell H1 (bold, centered, bordered) then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I0 / IF columns, rows 2-19.
$data = @(
    @(1, 4),
    @(1, 6),
    @(1, 4),
    @(1, 4),
    @(8, 8),
    @(1, 5),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(3, 3),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
